$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix header style A1: fill goes from red to white (theme "Background 1") ---
# The red fill made red header text invisible against a red background; switch the
# fill to the theme's white background color so the red text is legible again.
$ws.Range("A1").Interior.ThemeColor = 2
$ws.Range("A1").Interior.TintAndShade = 0

# --- 2. Rows 41-42 and 581-582 pick up the standard striped-data style (fontId 0 / fillId 2) ---
# Re-use the format already present on A43 (which already carries that style) via
# copy / paste-special so we don't mint new style records.
$ws.Range("A43").Copy()
$ws.Range("A41:A42").PasteSpecial(-4122)
$ws.Range("A43").Copy()
$ws.Range("A581:A582").PasteSpecial(-4122)

# --- 3. Column A width: 25.5 characters, matching the rest of the cleaned-up sheet ---
$ws.Columns("A").ColumnWidth = 24.6

# --- 4. Append the new ids (rows 583-674) ---
$values = @(12900612745,12900601934,12900586917,12900566978,12900543510,12900511357,12900463211,12900420036,12900369194,12900190549,12900169901,12900149072,12900114171,12899213448,12899134862,12898964027,12898642970,12898365222,12898343328,12898343247,12898343276,12898343271,12898343265,12898343258,12898343245,12898343252,12898343246,12898343257,12898343275,12898343249,12898343260,12898343250,12898343283,12898343264,12898343337,12898346382,12898330068,12871365441,12871365340,12871365208,12870703808,12870608853,12870608644,12870608540,12870608412,12870605393,12870605319,12870605251,12870605086,12870605168,12870608901,12870602990,12870570962,12870563163,12870562919,12870562798,12870562628,12870562509,12870575866,12870535188,12870534563,12870536330,12870535776,12870554968,12870534440,12870487929,12870538153,12870538027,12870537898,12870537764,12870537662,12870487102,12870486428,12870487774,12870564403,12870532864,12870543836,12870506467,12870484493,12870514486,12870514132,12870489185,12870483349,12870460246,12741081295,12720598465,12720595503,12676540645,12568164440,12707383608,12706315010,12869063756)
$startRow = 583
$row = $startRow
foreach ($v in $values) {
  $ws.Cells.Item($row, 1).Value = $v
  $row = $row + 1
}
$endRow = $row - 1

# Apply the same striped-data style to the new block, except the very last row.
$ws.Range("A43").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)

# The final row (674) is left unshaded - explicitly "no fill" rather than inheriting
# the striped background.
$lastCell = $ws.Cells.Item($endRow, 1)
$lastCell.Interior.ColorIndex = -4142

# --- 5. Put the selection on the new last cell, like the original author left it ---
$ws.Cells.Item($endRow, 1).Select()
